$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column C for rows 2-36
# from serial date 45680 (2025-01-23) to 45681 (2025-01-24)
for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 3).Value = 45681
}
